$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates for the two pairs of rows (weekly data correction):
# D4, D5: 2022-09-28 (44832) -> 2022-10-04 (44838)
# D6, D7: 2022-10-04 (44838) -> 2022-09-28 (44832)
$ws.Range("D4").Value = 44838
$ws.Range("D5").Value = 44838
$ws.Range("D6").Value = 44832
$ws.Range("D7").Value = 44832
